$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# B1, D1, E1 keep the same text; only C1's header text changes.
$ws.Range("C1").Value = "Execution time (in ms)"
$ws.Range("F1").Value = "Max frequency"
$ws.Range("G1").Value = "Number of different frequency"

# --- Row 2: MGM / scen02 ---
# C2 switches from a formatted text string ("28,718 ms") to a plain number.
$ws.Range("C2").Value = 28718
$ws.Range("F2").Value = 792
$ws.Range("G2").Value = 44
# Row 2's numeric cells use the thousands-separator number format (style index 1).
$ws.Range("C2:C2").NumberFormat = "#,##0"
$ws.Range("F2:G2").NumberFormat = "#,##0"

# --- Row 3: DSA / scen02 ---
$ws.Range("C3").Value = 88488
$ws.Range("F3").Value = 792
$ws.Range("G3").Value = 44

# --- Row 4 (new): MGM / scen05 ---
$ws.Range("A4").Value = "MGM"
$ws.Range("B4").Value = "scen05"
$ws.Range("C4").Value = 25250
$ws.Range("D4").Value = 2078400
$ws.Range("E4").Value = 43627015
$ws.Range("F4").Value = 792
$ws.Range("G4").Value = 44

# --- Row 5 (new): DSA / scen05 ---
$ws.Range("A5").Value = "DSA"
$ws.Range("B5").Value = "scen05"
$ws.Range("C5").Value = 26341
$ws.Range("D5").Value = 1034004
$ws.Range("E5").Value = 21771879
$ws.Range("F5").Value = 792
$ws.Range("G5").Value = 44

# --- New column widths (target widths are 20 and 27.88671875 characters;
# the COM ColumnWidth setter here quantizes to whole-pixel steps, so these
# inputs are chosen to land on the closest achievable stored width) ---
$ws.Columns("F").ColumnWidth = 19.15
$ws.Columns("G").ColumnWidth = 27

# --- Selection moves to G6 ---
$ws.Range("G6").Select()
